$d = $word.ActiveDocument

# --- Locate "May" inside "(Expected Graduation: May 2016)" -------------------------
# Search for the unique surrounding phrase so we don't disturb the many other
# occurrences of "May" used elsewhere in the resume (month abbreviations, etc.).
$context = $d.Content
$found = $context.Find.Execute("Expected Graduation: May 2016", $true, $false, $false, `
                                $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Expected Graduation: May 2016' in the document."
}

$contextText = $context.Text
$mayOffsetInContext = $contextText.IndexOf("May")
$mayStart = $context.Start + $mayOffsetInContext
$mayEnd = $mayStart + 3          # just the 3 characters "May"

# --- Replace "May" with "December", keeping that run's own formatting --------------
$mayRange = $d.Range($mayStart, $mayEnd)
$mayRange.Text = "December"

$decEnd = $mayStart + 8          # length of "December"

# The engine coalesces the (formatting-identical) runs that follow the edit point
# (" ", "2016", ")") into a single run. Re-split them back into their original three
# runs by round-tripping a character property (Bold off->on->off) on each piece; this
# forces the engine to re-materialize run boundaries there without altering the
# run's actual appearance/formatting (the "cs" font etc. stay intact, unlike
# reassigning Font.Name/NameAscii, which strips the w:cs attribute entirely).
$spaceStart = $decEnd
$spaceEnd = $spaceStart + 1
$yearStart = $spaceEnd
$yearEnd = $yearStart + 4
$parenStart = $yearEnd
$parenEnd = $parenStart + 1

foreach ($bounds in @(@($yearStart, $yearEnd), @($parenStart, $parenEnd))) {
    $piece = $d.Range($bounds[0], $bounds[1])
    $piece.Font.Bold = 1
    $piece.Font.Bold = 0
}

# --- Move the "_GoBack" bookmark so it again sits right after the edited word ------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$newBookmarkPoint = $d.Range($decEnd, $decEnd)
$d.Bookmarks.Add("_GoBack", $newBookmarkPoint)
